# Rename the "Forecast" threshold labels (column A) to "Pipeline" across the
# whole used range, then restore the cursor/selection to J31 with the view
# scrolled back to the top-left (matching the post-edit sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every occurrence of "Forecast" with "Pipeline" in the label column.
# (Covers: "Facturacion por Forecast >= 60%", "Facturacion por Forecast< 60%",
#  "Facturacion por Forecast < 60%", "Costo de venta por Forecast >= 60%",
#  "Costo de venta por Forecast < 60%".)
$rng = $ws.Range("A1:A85")
[void]$rng.Replace("Forecast", "Pipeline")

# Move/restore the active selection to J31 (also resets the scrolled
# top-left cell back to the sheet's natural origin).
[void]$ws.Range("J31").Select()
